$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("YDS")

$ws.Range("B2").Value = "NIG(0.9009464725209637, 0.6455690948284724, 1.1889830929412324, 3.1678282639942745)"
$ws.Range("C2").Value = "JSU(-0.8692851019365995, 1.052748650101893, 4.384017887912757, 4.421979114585777)"
$ws.Range("D2").Value = "NIG(0.5501148123438273, 0.2850034474612263, 2.7407730556249605, 2.855733732087514)"
$ws.Range("E2").Value = "JSU(-1.0436728175364258, 1.0843471331452244, 3.4239441456569133, 4.8532444915904875)"
